$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.67
$ws.Range("G2").Value = 1.8
$ws.Range("J2").Value = 3.9
$ws.Range("N2").Value = 3.75
$ws.Range("O2").Value = 1.29
$ws.Range("Q2").Value = 1.87
$ws.Range("T2").Value = 1.84
$ws.Range("U2").Value = 1.98
$ws.Range("W2").Value = 2.26
$ws.Range("AB2").Value = 12.5
$ws.Range("AF2").Value = 23
$ws.Range("AG2").Value = 19
$ws.Range("AJ2").Value = 900
$ws.Range("AK2").Value = 70
$ws.Range("AN2").Value = 28
$ws.Range("F3").Value = 1.72
$ws.Range("I3").Value = 5.9
$ws.Range("L3").Value = 1.33
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 21
$ws.Range("AB3").Value = 11
$ws.Range("AD3").Value = 20
$ws.Range("AE3").Value = 1000
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 300
$ws.Range("AJ3").Value = 20
$ws.Range("AK3").Value = 18.5
$ws.Range("AL3").Value = 90
$ws.Range("AM3").Value = 580
$ws.Range("AN3").Value = 9.800000000000001
$ws.Range("F5").Value = 1.39
$ws.Range("G5").Value = 1.41
$ws.Range("H5").Value = 9.6
$ws.Range("J5").Value = 5.3
$ws.Range("K5").Value = 5.5
$ws.Range("P5").Value = 2.1
$ws.Range("Q5").Value = 1.85
$ws.Range("V5").Value = 1.11
$ws.Range("W5").Value = 3.4
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 50
$ws.Range("Z5").Value = 90
$ws.Range("AA5").Value = 390
$ws.Range("AH5").Value = 30
$ws.Range("AI5").Value = 160
$ws.Range("AL5").Value = 42
$ws.Range("AO5").Value = 250
$ws.Range("G6").Value = 4.8
$ws.Range("H6").Value = 2.04
$ws.Range("I6").Value = 2.1
$ws.Range("J6").Value = 3.35
$ws.Range("K6").Value = 3.45
$ws.Range("N6").Value = 3
$ws.Range("P6").Value = 1.68
$ws.Range("Q6").Value = 2.38
$ws.Range("R6").Value = 1.25
$ws.Range("T6").Value = 2.06
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 1.27
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 7.4
$ws.Range("Z6").Value = 11.5
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 13
$ws.Range("AE6").Value = 26
$ws.Range("AF6").Value = 30
$ws.Range("AG6").Value = 19
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 110
$ws.Range("AK6").Value = 70
$ws.Range("AL6").Value = 90
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 100
$ws.Range("AO6").Value = 23
$ws.Range("G7").Value = 2.32
$ws.Range("H7").Value = 3.65
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 3.2
$ws.Range("P7").Value = 1.73
$ws.Range("R7").Value = 1.27
$ws.Range("T7").Value = 1.9
$ws.Range("U7").Value = 1.95
$ws.Range("W7").Value = 1.75
$ws.Range("X7").Value = 11.5
$ws.Range("Z7").Value = 26
$ws.Range("AA7").Value = 85
$ws.Range("AC7").Value = 7.4
$ws.Range("AD7").Value = 15.5
$ws.Range("AE7").Value = 55
$ws.Range("AF7").Value = 13.5
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 20
$ws.Range("AI7").Value = 75
$ws.Range("AJ7").Value = 30
$ws.Range("AK7").Value = 27
$ws.Range("AN7").Value = 24
$ws.Range("AO7").Value = 65
$ws.Range("F8").Value = 2.54
$ws.Range("G8").Value = 2.6
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 3.35
$ws.Range("O8").Value = 1.48
$ws.Range("P8").Value = 1.66
$ws.Range("S8").Value = 4.8
$ws.Range("V8").Value = 1.42
$ws.Range("W8").Value = 1.62
$ws.Range("AE8").Value = 46
$ws.Range("AG8").Value = 12
$ws.Range("AH8").Value = 21
$ws.Range("AI8").Value = 65
$ws.Range("AJ8").Value = 38
$ws.Range("AL8").Value = 55
$ws.Range("AN8").Value = 38
$ws.Range("G9").Value = 2.7
$ws.Range("H9").Value = 3.15
$ws.Range("L9").Value = 1.53
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 2.88
$ws.Range("P9").Value = 1.62
$ws.Range("Q9").Value = 2.48
$ws.Range("T9").Value = 2.02
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 25
$ws.Range("AA9").Value = 1000
$ws.Range("AB9").Value = 10.5
$ws.Range("AC9").Value = 7.6
$ws.Range("AD9").Value = 17.5
$ws.Range("AE9").Value = 55
$ws.Range("AG9").Value = 13
$ws.Range("AH9").Value = 980
$ws.Range("AI9").Value = 1000
$ws.Range("AJ9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AM9").Value = 190
$ws.Range("AN9").Value = 1000
$ws.Range("G10").Value = 2.44
$ws.Range("H10").Value = 3.7
$ws.Range("L10").Value = 1.6
$ws.Range("N10").Value = 2.62
$ws.Range("Q10").Value = 2.74
$ws.Range("T10").Value = 2.2
$ws.Range("U10").Value = 1.77
$ws.Range("V10").Value = 1.35
$ws.Range("W10").Value = 1.69
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 28
$ws.Range("AA10").Value = 990
$ws.Range("AC10").Value = 7.2
$ws.Range("AD10").Value = 22
$ws.Range("AE10").Value = 80
$ws.Range("AF10").Value = 15
$ws.Range("AG10").Value = 12.5
$ws.Range("AH10").Value = 38
$ws.Range("AI10").Value = 110
$ws.Range("AJ10").Value = 40
$ws.Range("AK10").Value = 40
$ws.Range("AN10").Value = 65
$ws.Range("AO10").Value = 120
$ws.Range("G11").Value = 2.4
$ws.Range("I11").Value = 3.55
$ws.Range("N11").Value = 3.25
$ws.Range("S11").Value = 4.2
$ws.Range("T11").Value = 1.87
$ws.Range("W11").Value = 1.71
$ws.Range("X11").Value = 12
$ws.Range("Y11").Value = 12
$ws.Range("Z11").Value = 1000
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 9
$ws.Range("AC11").Value = 7.8
$ws.Range("AE11").Value = 980
$ws.Range("AG11").Value = 11.5
$ws.Range("AH11").Value = 980
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 1000
$ws.Range("AK11").Value = 1000
$ws.Range("AL11").Value = 1000
$ws.Range("AN11").Value = 1000
$ws.Range("F12").Value = 3.2
$ws.Range("G12").Value = 3.3
$ws.Range("H12").Value = 2.62
$ws.Range("I12").Value = 2.68
$ws.Range("N12").Value = 2.68
$ws.Range("O12").Value = 1.55
$ws.Range("U12").Value = 1.81
$ws.Range("V12").Value = 1.59
$ws.Range("W12").Value = 1.43
$ws.Range("X12").Value = 8.800000000000001
$ws.Range("Y12").Value = 8
$ws.Range("Z12").Value = 15.5
$ws.Range("AA12").Value = 42
$ws.Range("AB12").Value = 9.199999999999999
$ws.Range("AC12").Value = 7.2
$ws.Range("AD12").Value = 13
$ws.Range("AE12").Value = 38
$ws.Range("AF12").Value = 19.5
$ws.Range("AH12").Value = 24
$ws.Range("AI12").Value = 65
$ws.Range("AJ12").Value = 60
$ws.Range("AK12").Value = 50
$ws.Range("AL12").Value = 75
$ws.Range("AM12").Value = 180
$ws.Range("F13").Value = 2.08
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 4.4
$ws.Range("J13").Value = 2.96
$ws.Range("K13").Value = 4.4
$ws.Range("L13").Value = 1.35
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 3
$ws.Range("O13").Value = 1.33
$ws.Range("P13").Value = 1.8
$ws.Range("Q13").Value = 1.94
$ws.Range("R13").Value = 1.31
$ws.Range("S13").Value = 3.1
$ws.Range("T13").Value = 1.74
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.29
$ws.Range("W13").Value = 1.64
$ws.Range("AC13").Value = 980
$ws.Range("F14").Value = 2.16
$ws.Range("G14").Value = 2.36
$ws.Range("H14").Value = 3.35
$ws.Range("J14").Value = 3.3
$ws.Range("N14").Value = 3
$ws.Range("Q14").Value = 2.16
$ws.Range("V14").Value = 1.34
$ws.Range("W14").Value = 1.73
$ws.Range("X14").Value = 14
$ws.Range("Y14").Value = 14
$ws.Range("Z14").Value = 32
$ws.Range("AA14").Value = 1000
$ws.Range("AB14").Value = 10
$ws.Range("AC14").Value = 9.6
$ws.Range("AD14").Value = 18.5
$ws.Range("AE14").Value = 1000
$ws.Range("AF14").Value = 16.5
$ws.Range("AG14").Value = 13
$ws.Range("AH14").Value = 24
$ws.Range("AI14").Value = 1000
$ws.Range("AJ14").Value = 1000
$ws.Range("AK14").Value = 34
$ws.Range("AL14").Value = 60
$ws.Range("AO14").Value = 75
